$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MACRO_SCORE (column N) recalculated for all four data rows.
$ws.Range("N2:N5").Value = 54.84087454262382
